$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-21 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-22 Saturday", 2) | Out-Null
$d.Content.Find.Execute("470÷6=78, 2", $true, $false, $false, $false, $false, $true, 1, $false, "460÷2=230, 0", 2) | Out-Null
$d.Content.Find.Execute("569÷9=63, 2", $true, $false, $false, $false, $false, $true, 1, $false, "642÷6=107, 0", 2) | Out-Null
$d.Content.Find.Execute("142÷5=28, 2", $true, $false, $false, $false, $false, $true, 1, $false, "788÷3=262, 2", 2) | Out-Null
$d.Content.Find.Execute("424÷7=60, 4", $true, $false, $false, $false, $false, $true, 1, $false, "114÷9=12, 6", 2) | Out-Null
$d.Content.Find.Execute("656÷2=328, 0", $true, $false, $false, $false, $false, $true, 1, $false, "432÷3=144, 0", 2) | Out-Null
$d.Content.Find.Execute("151÷8=18, 7", $true, $false, $false, $false, $false, $true, 1, $false, "307÷8=38, 3", 2) | Out-Null
$d.Content.Find.Execute("555÷5=111, 0", $true, $false, $false, $false, $false, $true, 1, $false, "493÷3=164, 1", 2) | Out-Null
$d.Content.Find.Execute("880÷8=110, 0", $true, $false, $false, $false, $false, $true, 1, $false, "469÷3=156, 1", 2) | Out-Null
$d.Content.Find.Execute("371÷7=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "696÷5=139, 1", 2) | Out-Null
$d.Content.Find.Execute("358÷6=59, 4", $true, $false, $false, $false, $false, $true, 1, $false, "176÷4=44, 0", 2) | Out-Null
$d.Content.Find.Execute("139÷3=46, 1", $true, $false, $false, $false, $false, $true, 1, $false, "388÷9=43, 1", 2) | Out-Null
$d.Content.Find.Execute("235÷2=117, 1", $true, $false, $false, $false, $false, $true, 1, $false, "467÷7=66, 5", 2) | Out-Null
$d.Content.Find.Execute("936÷8=117, 0", $true, $false, $false, $false, $false, $true, 1, $false, "799÷8=99, 7", 2) | Out-Null
$d.Content.Find.Execute("729÷2=364, 1", $true, $false, $false, $false, $false, $true, 1, $false, "836÷3=278, 2", 2) | Out-Null
$d.Content.Find.Execute("549÷6=91, 3", $true, $false, $false, $false, $false, $true, 1, $false, "351÷2=175, 1", 2) | Out-Null
$d.Content.Find.Execute("223÷6=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "224÷5=44, 4", 2) | Out-Null
$d.Content.Find.Execute("441÷8=55, 1", $true, $false, $false, $false, $false, $true, 1, $false, "298÷6=49, 4", 2) | Out-Null
$d.Content.Find.Execute("817÷5=163, 2", $true, $false, $false, $false, $false, $true, 1, $false, "769÷4=192, 1", 2) | Out-Null
$d.Content.Find.Execute("990÷2=495, 0", $true, $false, $false, $false, $false, $true, 1, $false, "293÷9=32, 5", 2) | Out-Null
$d.Content.Find.Execute("419÷5=83, 4", $true, $false, $false, $false, $false, $true, 1, $false, "116÷2=58, 0", 2) | Out-Null
$d.Content.Find.Execute("513÷3=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "642÷8=80, 2", 2) | Out-Null
$d.Content.Find.Execute("636÷3=212, 0", $true, $false, $false, $false, $false, $true, 1, $false, "194÷3=64, 2", 2) | Out-Null
$d.Content.Find.Execute("637÷6=106, 1", $true, $false, $false, $false, $false, $true, 1, $false, "227÷6=37, 5", 2) | Out-Null
$d.Content.Find.Execute("375÷5=75, 0", $true, $false, $false, $false, $false, $true, 1, $false, "973÷4=243, 1", 2) | Out-Null
$d.Content.Find.Execute("142÷7=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "454÷5=90, 4", 2) | Out-Null
